# Auto-generated Excel COM-interop script to update Cactuar Profits data
# Applies per-row H:N numeric updates across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 970.1667  # H12: 816.7273 -> 970.1667
$ws.Cells.Item(12, 9).Value = 961.5  # I12: 874.5 -> 961.5
$ws.Cells.Item(12, 10).Value = 987.5  # J12: 747.4 -> 987.5
$ws.Cells.Item(12, 11).Value = 961.5  # K12: 874.5 -> 961.5
$ws.Cells.Item(12, 12).Value = 987.5  # L12: 747.4 -> 987.5
$ws.Cells.Item(12, 13).Value = -791.5  # M12: -704.5 -> -791.5
$ws.Cells.Item(12, 14).Value = -1327.5  # N12: -1087.4 -> -1327.5

$ws.Cells.Item(38, 8).Value = 4286.7  # H38: 4110.905 -> 4286.7
$ws.Cells.Item(38, 9).Value = 1794.7273  # I38: 1694.75 -> 1794.7273
$ws.Cells.Item(38, 11).Value = 5384.1819  # K38: 5084.25 -> 5384.1819
$ws.Cells.Item(38, 13).Value = -5012.1819  # M38: -4712.25 -> -5012.1819

$ws.Cells.Item(61, 8).Value = 338.66666  # H61: 141.66667 -> 338.66666
$ws.Cells.Item(61, 9).Value = 338.66666  # I61: 141.66667 -> 338.66666
$ws.Cells.Item(61, 11).Value = 1015.99998  # K61: 425.00001 -> 1015.99998
$ws.Cells.Item(61, 13).Value = -843.9999799999999  # M61: -253.00001 -> -843.9999799999999

$ws.Cells.Item(64, 8).Value = 10872414  # H64: 6214393 -> 10872414
$ws.Cells.Item(64, 10).Value = 3799  # J64: 3748.5 -> 3799
$ws.Cells.Item(64, 12).Value = 3799  # L64: 3748.5 -> 3799
$ws.Cells.Item(64, 14).Value = -4295  # N64: -4244.5 -> -4295

$ws.Cells.Item(67, 8).Value = 10872414  # H67: 6214393 -> 10872414
$ws.Cells.Item(67, 10).Value = 3799  # J67: 3748.5 -> 3799
$ws.Cells.Item(67, 12).Value = 3799  # L67: 3748.5 -> 3799
$ws.Cells.Item(67, 14).Value = -5515  # N67: -5464.5 -> -5515

$ws.Cells.Item(69, 8).Value = 10000  # H69: 0 -> 10000
$ws.Cells.Item(69, 10).Value = 10000  # J69: 0 -> 10000
$ws.Cells.Item(69, 12).Value = 30000  # L69: 0 -> 30000
$ws.Cells.Item(69, 14).Value = -31748  # N69: None -> -31748

$ws.Cells.Item(72, 8).Value = 10000  # H72: 0 -> 10000
$ws.Cells.Item(72, 10).Value = 10000  # J72: 0 -> 10000
$ws.Cells.Item(72, 12).Value = 90000  # L72: 0 -> 90000
$ws.Cells.Item(72, 14).Value = -98736  # N72: None -> -98736

$ws.Cells.Item(100, 8).Value = 819.6429000000001  # H100: 791.3333 -> 819.6429000000001
$ws.Cells.Item(100, 9).Value = 813.46155  # I100: 783.5714 -> 813.46155
$ws.Cells.Item(100, 11).Value = 813.46155  # K100: 783.5714 -> 813.46155
$ws.Cells.Item(100, 13).Value = -272.46155  # M100: -242.5714 -> -272.46155

$ws.Cells.Item(106, 8).Value = 19610106  # H106: 17546072 -> 19610106
$ws.Cells.Item(106, 9).Value = 22224000  # I106: 19609622 -> 22224000
$ws.Cells.Item(106, 11).Value = 22224000  # K106: 19609622 -> 22224000
$ws.Cells.Item(106, 13).Value = -22223369  # M106: -19608991 -> -22223369

$ws.Cells.Item(121, 8).Value = 4829.6665  # H121: 4808.357 -> 4829.6665
$ws.Cells.Item(121, 10).Value = 4829.6665  # J121: 4808.357 -> 4829.6665
$ws.Cells.Item(121, 12).Value = 14488.9995  # L121: 14425.071 -> 14488.9995
$ws.Cells.Item(121, 14).Value = -17982.9995  # N121: -17919.071 -> -17982.9995

$ws.Cells.Item(132, 8).Value = 16125.269  # H132: 16508.4 -> 16125.269
$ws.Cells.Item(132, 9).Value = 3707.7693  # I132: 3950.0833 -> 3707.7693
$ws.Cells.Item(132, 11).Value = 11123.3079  # K132: 11850.2499 -> 11123.3079
$ws.Cells.Item(132, 13).Value = -8593.3079  # M132: -9320.249899999999 -> -8593.3079

$ws.Cells.Item(135, 8).Value = 2626.1936  # H135: 2545.5312 -> 2626.1936
$ws.Cells.Item(135, 9).Value = 814.9545000000001  # I135: 781.4783 -> 814.9545000000001
$ws.Cells.Item(135, 11).Value = 7334.5905  # K135: 7033.3047 -> 7334.5905
$ws.Cells.Item(135, 13).Value = -4799.5905  # M135: -4498.3047 -> -4799.5905

$ws.Cells.Item(138, 8).Value = 3543.4243  # H138: 3613.2222 -> 3543.4243
$ws.Cells.Item(138, 10).Value = 4260.28  # J138: 4352.4136 -> 4260.28
$ws.Cells.Item(138, 12).Value = 12780.84  # L138: 13057.2408 -> 12780.84
$ws.Cells.Item(138, 14).Value = -23060.84  # N138: -23337.2408 -> -23060.84

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11489.926  # H32: 11917.641 -> 11489.926
$ws.Cells.Item(32, 9).Value = 11117.209  # I32: 11283.97 -> 11117.209
$ws.Cells.Item(32, 10).Value = 13273.643  # J32: 15402.833 -> 13273.643
$ws.Cells.Item(32, 11).Value = 11117.209  # K32: 11283.97 -> 11117.209
$ws.Cells.Item(32, 12).Value = 13273.643  # L32: 15402.833 -> 13273.643
$ws.Cells.Item(32, 13).Value = -10830.209  # M32: -10996.97 -> -10830.209
$ws.Cells.Item(32, 14).Value = -13847.643  # N32: -15976.833 -> -13847.643

$ws.Cells.Item(45, 8).Value = 3963.5715  # H45: 3280.625 -> 3963.5715
$ws.Cells.Item(45, 9).Value = 2000  # I45: 1833 -> 2000
$ws.Cells.Item(45, 10).Value = 4290.8335  # J45: 4149.2 -> 4290.8335
$ws.Cells.Item(45, 11).Value = 2000  # K45: 1833 -> 2000
$ws.Cells.Item(45, 12).Value = 4290.8335  # L45: 4149.2 -> 4290.8335
$ws.Cells.Item(45, 13).Value = -1623  # M45: -1456 -> -1623
$ws.Cells.Item(45, 14).Value = -5044.8335  # N45: -4903.2 -> -5044.8335

$ws.Cells.Item(97, 8).Value = 308.1154  # H97: 315.5 -> 308.1154
$ws.Cells.Item(97, 9).Value = 312.5238  # I97: 321.66666 -> 312.5238
$ws.Cells.Item(97, 11).Value = 312.5238  # K97: 321.66666 -> 312.5238
$ws.Cells.Item(97, 13).Value = 183.4762  # M97: 174.33334 -> 183.4762

$ws.Cells.Item(102, 8).Value = 597404.25  # H102: 490842.25 -> 597404.25
$ws.Cells.Item(102, 9).Value = 686415.4  # I102: 549263.75 -> 686415.4
$ws.Cells.Item(102, 11).Value = 686415.4  # K102: 549263.75 -> 686415.4
$ws.Cells.Item(102, 13).Value = -684793.4  # M102: -547641.75 -> -684793.4

$ws.Cells.Item(109, 8).Value = 95000  # H109: 0 -> 95000
$ws.Cells.Item(109, 10).Value = 95000  # J109: 0 -> 95000
$ws.Cells.Item(109, 12).Value = 95000  # L109: 0 -> 95000
$ws.Cells.Item(109, 14).Value = -97774  # N109: None -> -97774

$ws.Cells.Item(122, 8).Value = 5970.1055  # H122: 5506.2856 -> 5970.1055
$ws.Cells.Item(122, 9).Value = 3073.7273  # I122: 2770.077 -> 3073.7273
$ws.Cells.Item(122, 11).Value = 9221.1819  # K122: 8310.231 -> 9221.1819
$ws.Cells.Item(122, 13).Value = -6771.1819  # M122: -5860.231 -> -6771.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 500  # H11: 352 -> 500
$ws.Cells.Item(11, 9).Value = 500  # I11: 352 -> 500
$ws.Cells.Item(11, 11).Value = 500  # K11: 352 -> 500
$ws.Cells.Item(11, 13).Value = -360  # M11: -212 -> -360

$ws.Cells.Item(20, 8).Value = 2838.6904  # H20: 2887.439 -> 2838.6904
$ws.Cells.Item(20, 9).Value = 2290.276  # I20: 2342.0715 -> 2290.276
$ws.Cells.Item(20, 11).Value = 2290.276  # K20: 2342.0715 -> 2290.276
$ws.Cells.Item(20, 13).Value = -2043.276  # M20: -2095.0715 -> -2043.276

$ws.Cells.Item(99, 8).Value = 1489533.9  # H99: 1303466.8 -> 1489533.9
$ws.Cells.Item(99, 9).Value = 1737456.1  # I99: 1489390.5 -> 1737456.1
$ws.Cells.Item(99, 11).Value = 1737456.1  # K99: 1489390.5 -> 1737456.1
$ws.Cells.Item(99, 13).Value = -1735958.1  # M99: -1487892.5 -> -1735958.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5788.8213  # H31: 5896.327 -> 5788.8213
$ws.Cells.Item(31, 9).Value = 2353.2942  # I31: 2492.1333 -> 2353.2942
$ws.Cells.Item(31, 10).Value = 7286.359  # J31: 7172.9 -> 7286.359
$ws.Cells.Item(31, 11).Value = 2353.2942  # K31: 2492.1333 -> 2353.2942
$ws.Cells.Item(31, 12).Value = 7286.359  # L31: 7172.9 -> 7286.359
$ws.Cells.Item(31, 13).Value = -2058.2942  # M31: -2197.1333 -> -2058.2942
$ws.Cells.Item(31, 14).Value = -7876.359  # N31: -7762.9 -> -7876.359

$ws.Cells.Item(34, 8).Value = 5788.8213  # H34: 5896.327 -> 5788.8213
$ws.Cells.Item(34, 9).Value = 2353.2942  # I34: 2492.1333 -> 2353.2942
$ws.Cells.Item(34, 10).Value = 7286.359  # J34: 7172.9 -> 7286.359
$ws.Cells.Item(34, 11).Value = 2353.2942  # K34: 2492.1333 -> 2353.2942
$ws.Cells.Item(34, 12).Value = 7286.359  # L34: 7172.9 -> 7286.359
$ws.Cells.Item(34, 13).Value = -2151.2942  # M34: -2290.1333 -> -2151.2942
$ws.Cells.Item(34, 14).Value = -7690.359  # N34: -7576.9 -> -7690.359

$ws.Cells.Item(99, 8).Value = 8180.4116  # H99: 7095.2383 -> 8180.4116
$ws.Cells.Item(99, 10).Value = 8759.076999999999  # J99: 7282.4116 -> 8759.076999999999
$ws.Cells.Item(99, 12).Value = 8759.076999999999  # L99: 7282.4116 -> 8759.076999999999
$ws.Cells.Item(99, 14).Value = -11755.077  # N99: -10278.4116 -> -11755.077

$ws.Cells.Item(105, 8).Value = 3247774.5  # H105: 2273581.8 -> 3247774.5
$ws.Cells.Item(105, 9).Value = 4546084  # I105: 2841477.2 -> 4546084
$ws.Cells.Item(105, 11).Value = 4546084  # K105: 2841477.2 -> 4546084
$ws.Cells.Item(105, 13).Value = -4544337  # M105: -2839730.2 -> -4544337

$ws.Cells.Item(107, 8).Value = 1818750.8  # H107: 2020834.1 -> 1818750.8
$ws.Cells.Item(107, 9).Value = 3636883.5  # I107: 4546104 -> 3636883.5
$ws.Cells.Item(107, 11).Value = 3636883.5  # K107: 4546104 -> 3636883.5
$ws.Cells.Item(107, 13).Value = -3634963.5  # M107: -4544184 -> -3634963.5

$ws.Cells.Item(109, 8).Value = 80000  # H109: 79997 -> 80000
$ws.Cells.Item(109, 10).Value = 80000  # J109: 79997 -> 80000
$ws.Cells.Item(109, 12).Value = 80000  # L109: 79997 -> 80000
$ws.Cells.Item(109, 14).Value = -82080  # N109: -82077 -> -82080

$ws.Cells.Item(126, 8).Value = 8180.4116  # H126: 7095.2383 -> 8180.4116
$ws.Cells.Item(126, 10).Value = 8759.076999999999  # J126: 7282.4116 -> 8759.076999999999
$ws.Cells.Item(126, 12).Value = 26277.231  # L126: 21847.2348 -> 26277.231
$ws.Cells.Item(126, 14).Value = -31217.231  # N126: -26787.2348 -> -31217.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 11099.667  # H3: 13333.333 -> 11099.667
$ws.Cells.Item(3, 9).Value = 799  # I3: 0 -> 799
$ws.Cells.Item(3, 10).Value = 16250  # J3: 13333.333 -> 16250
$ws.Cells.Item(3, 11).Value = 2397  # K3: 0 -> 2397
$ws.Cells.Item(3, 12).Value = 48750  # L3: 39999.999 -> 48750
$ws.Cells.Item(3, 13).Value = -2285  # M3: None -> -2285
$ws.Cells.Item(3, 14).Value = -48974  # N3: -40223.999 -> -48974

$ws.Cells.Item(55, 8).Value = 2707.3635  # H55: 2654.8696 -> 2707.3635
$ws.Cells.Item(55, 10).Value = 3500.1875  # J55: 3382.5293 -> 3500.1875
$ws.Cells.Item(55, 12).Value = 10500.5625  # L55: 10147.5879 -> 10500.5625
$ws.Cells.Item(55, 14).Value = -10854.5625  # N55: -10501.5879 -> -10854.5625

$ws.Cells.Item(131, 8).Value = 14150808  # H131: 14638698 -> 14150808
$ws.Cells.Item(131, 10).Value = 14638732  # J131: 15161472 -> 14638732
$ws.Cells.Item(131, 12).Value = 43916196  # L131: 45484416 -> 43916196
$ws.Cells.Item(131, 14).Value = -43926276  # N131: -45494496 -> -43926276

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(55, 8).Value = 0  # H55: 2999 -> 0
$ws.Cells.Item(55, 9).Value = 0  # I55: 2999 -> 0
$ws.Cells.Item(55, 11).Value = 0  # K55: 2999 -> 0
$ws.Cells.Item(55, 13).ClearContents()  # M55: -2672 -> (removed)

$ws.Cells.Item(70, 8).Value = 3502280  # H70: 3793761.8 -> 3502280
$ws.Cells.Item(70, 9).Value = 4550664  # I70: 5055793.5 -> 4550664
$ws.Cells.Item(70, 11).Value = 4550664  # K70: 5055793.5 -> 4550664
$ws.Cells.Item(70, 13).Value = -4550394  # M70: -5055523.5 -> -4550394

$ws.Cells.Item(73, 8).Value = 3502280  # H73: 3793761.8 -> 3502280
$ws.Cells.Item(73, 9).Value = 4550664  # I73: 5055793.5 -> 4550664
$ws.Cells.Item(73, 11).Value = 4550664  # K73: 5055793.5 -> 4550664
$ws.Cells.Item(73, 13).Value = -4549728  # M73: -5054857.5 -> -4549728

$ws.Cells.Item(97, 8).Value = 702.05  # H97: 703.45 -> 702.05
$ws.Cells.Item(97, 9).Value = 568.73334  # I97: 570.6 -> 568.73334
$ws.Cells.Item(97, 11).Value = 568.73334  # K97: 570.6 -> 568.73334
$ws.Cells.Item(97, 13).Value = -72.73334  # M97: -74.60000000000002 -> -72.73334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5382.9443  # H7: 5587.75 -> 5382.9443
$ws.Cells.Item(7, 9).Value = 3089.4  # I7: 2925.625 -> 3089.4
$ws.Cells.Item(7, 11).Value = 3089.4  # K7: 2925.625 -> 3089.4
$ws.Cells.Item(7, 13).Value = -2977.4  # M7: -2813.625 -> -2977.4

$ws.Cells.Item(40, 8).Value = 2822.25  # H40: 3200.0625 -> 2822.25
$ws.Cells.Item(40, 9).Value = 2822.25  # I40: 3399.8667 -> 2822.25
$ws.Cells.Item(40, 10).Value = 0  # J40: 203 -> 0
$ws.Cells.Item(40, 11).Value = 2822.25  # K40: 3399.8667 -> 2822.25
$ws.Cells.Item(40, 12).Value = 0  # L40: 203 -> 0
$ws.Cells.Item(40, 13).Value = -2686.25  # M40: -3263.8667 -> -2686.25
$ws.Cells.Item(40, 14).ClearContents()  # N40: -475 -> (removed)

$ws.Cells.Item(61, 8).Value = 2099.5715  # H61: 3224.5 -> 2099.5715
$ws.Cells.Item(61, 9).Value = 782.8333  # I61: 966 -> 782.8333
$ws.Cells.Item(61, 11).Value = 782.8333  # K61: 966 -> 782.8333
$ws.Cells.Item(61, 13).Value = -580.8333  # M61: -764 -> -580.8333

$ws.Cells.Item(93, 8).Value = 1910.3334  # H93: 2052.4375 -> 1910.3334
$ws.Cells.Item(93, 9).Value = 1967.875  # I93: 2138.5 -> 1967.875
$ws.Cells.Item(93, 11).Value = 1967.875  # K93: 2138.5 -> 1967.875
$ws.Cells.Item(93, 13).Value = -719.875  # M93: -890.5 -> -719.875

$ws.Cells.Item(100, 8).Value = 10728.667  # H100: 9672.143 -> 10728.667
$ws.Cells.Item(100, 9).Value = 2075  # I100: 2494.3333 -> 2075
$ws.Cells.Item(100, 11).Value = 2075  # K100: 2494.3333 -> 2075
$ws.Cells.Item(100, 13).Value = -1534  # M100: -1953.3333 -> -1534

$ws.Cells.Item(113, 8).Value = 2099.5715  # H113: 3224.5 -> 2099.5715
$ws.Cells.Item(113, 9).Value = 782.8333  # I113: 966 -> 782.8333
$ws.Cells.Item(113, 11).Value = 782.8333  # K113: 966 -> 782.8333
$ws.Cells.Item(113, 13).Value = 1387.1667  # M113: 1204 -> 1387.1667

$ws.Cells.Item(122, 8).Value = 13416.889  # H122: 11842.454 -> 13416.889
$ws.Cells.Item(122, 9).Value = 4960.8  # I122: 4887.5 -> 4960.8
$ws.Cells.Item(122, 10).Value = 23987  # J122: 20188.4 -> 23987
$ws.Cells.Item(122, 11).Value = 14882.4  # K122: 14662.5 -> 14882.4
$ws.Cells.Item(122, 12).Value = 71961  # L122: 60565.2 -> 71961
$ws.Cells.Item(122, 13).Value = -12432.4  # M122: -12212.5 -> -12432.4
$ws.Cells.Item(122, 14).Value = -76861  # N122: -65465.2 -> -76861

$ws.Cells.Item(126, 8).Value = 5382.9443  # H126: 5587.75 -> 5382.9443
$ws.Cells.Item(126, 9).Value = 3089.4  # I126: 2925.625 -> 3089.4
$ws.Cells.Item(126, 11).Value = 9268.200000000001  # K126: 8776.875 -> 9268.200000000001
$ws.Cells.Item(126, 13).Value = -6798.200000000001  # M126: -6306.875 -> -6798.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 7654.2856  # H96: 7764.1665 -> 7654.2856
$ws.Cells.Item(96, 9).Value = 6989.6665  # I96: 6987 -> 6989.6665
$ws.Cells.Item(96, 11).Value = 6989.6665  # K96: 6987 -> 6989.6665
$ws.Cells.Item(96, 13).Value = -5616.6665  # M96: -5614 -> -5616.6665

$ws.Cells.Item(100, 8).Value = 1234141.8  # H100: 1016398.56 -> 1234141.8
$ws.Cells.Item(100, 9).Value = 1569513.9  # I100: 1233245.9 -> 1569513.9
$ws.Cells.Item(100, 11).Value = 3139027.8  # K100: 2466491.8 -> 3139027.8
$ws.Cells.Item(100, 13).Value = -3138486.8  # M100: -2465950.8 -> -3138486.8

$ws.Cells.Item(122, 8).Value = 3704.1  # H122: 3486.6 -> 3704.1
$ws.Cells.Item(122, 9).Value = 3782.1052  # I122: 3553.158 -> 3782.1052
$ws.Cells.Item(122, 11).Value = 11346.3156  # K122: 10659.474 -> 11346.3156
$ws.Cells.Item(122, 13).Value = -8896.3156  # M122: -8209.474 -> -8896.3156

$ws.Cells.Item(126, 8).Value = 3331  # H126: 5187.5 -> 3331
$ws.Cells.Item(126, 9).Value = 2806.8572  # I126: 4583.3335 -> 2806.8572
$ws.Cells.Item(126, 11).Value = 8420.571599999999  # K126: 13750.0005 -> 8420.571599999999
$ws.Cells.Item(126, 13).Value = -5950.571599999999  # M126: -11280.0005 -> -5950.571599999999

$ws.Cells.Item(127, 8).Value = 58997  # H127: 59000 -> 58997
$ws.Cells.Item(127, 10).Value = 58997  # J127: 59000 -> 58997
$ws.Cells.Item(127, 12).Value = 58997  # L127: 59000 -> 58997
$ws.Cells.Item(127, 14).Value = -68917  # N127: -68920 -> -68917

$ws.Cells.Item(132, 8).Value = 16129929  # H132: 17242326 -> 16129929
$ws.Cells.Item(132, 9).Value = 1014.1818  # I132: 1097.6 -> 1014.1818
$ws.Cells.Item(132, 11).Value = 3042.5454  # K132: 3292.8 -> 3042.5454
$ws.Cells.Item(132, 13).Value = -512.5454  # M132: -762.7999999999997 -> -512.5454
